# Update the header row labels as part of the data masking logic rework.
# Translate the Malay column headers to their English equivalents.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Student_Information")

$ws.Range("A1").Value = "Name"
$ws.Range("D1").Value = "Home Address"
$ws.Range("J1").Value = "Age"
$ws.Range("M1").Value = "Parent Salary (RM)"
$ws.Range("H1").Value = "Place of Birth"
